# Auto-generated edit script: updates cached Leve market-price values
# across multiple worksheet tabs (ALC, ARM, BSM, CRP, CUL, GSM, LTW).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1319.5
$ws.Range("I17").Value = 1296
$ws.Range("J17").Value = 1343
$ws.Range("K17").Value = 3888
$ws.Range("L17").Value = 4029
$ws.Range("M17").Value = -3720
$ws.Range("N17").Value = -4365
$ws.Range("H18").Value = 21665.334
$ws.Range("I18").Value = 19999
$ws.Range("J18").Value = 24998
$ws.Range("K18").Value = 19999
$ws.Range("L18").Value = 24998
$ws.Range("M18").Value = -19715
$ws.Range("N18").Value = -25566
$ws.Range("H32").Value = 990.3333
$ws.Range("I32").Value = 1625
$ws.Range("J32").Value = 809
$ws.Range("K32").Value = 1625
$ws.Range("L32").Value = 809
$ws.Range("M32").Value = -1299
$ws.Range("N32").Value = -1461
$ws.Range("H33").Value = 966.55554
$ws.Range("I33").Value = 366.5
$ws.Range("J33").Value = 2166.6667
$ws.Range("K33").Value = 366.5
$ws.Range("L33").Value = 2166.6667
$ws.Range("M33").Value = -137.5
$ws.Range("N33").Value = -2624.6667
$ws.Range("H40").Value = 4983.2856
$ws.Range("J40").Value = 5379.2
$ws.Range("L40").Value = 5379.2
$ws.Range("N40").Value = -5729.2
$ws.Range("H41").Value = 656.6429000000001
$ws.Range("I41").Value = 466.22223
$ws.Range("J41").Value = 999.4
$ws.Range("K41").Value = 466.22223
$ws.Range("L41").Value = 999.4
$ws.Range("M41").Value = -26.22223000000002
$ws.Range("N41").Value = -1879.4
$ws.Range("H43").Value = 14500
$ws.Range("J43").Value = 14500
$ws.Range("L43").Value = 14500
$ws.Range("N43").Value = -14638
$ws.Range("H137").Value = 12911.728
$ws.Range("I137").Value = 17842.428
$ws.Range("J137").Value = 4283
$ws.Range("K137").Value = 53527.284
$ws.Range("L137").Value = 12849
$ws.Range("M137").Value = -50977.284
$ws.Range("N137").Value = -17949
$ws.Range("H138").Value = 2444
$ws.Range("I138").Value = 1556
$ws.Range("J138").Value = 3406
$ws.Range("K138").Value = 4668
$ws.Range("L138").Value = 10218
$ws.Range("M138").Value = 472
$ws.Range("N138").Value = -20498
$ws.Range("H141").Value = 3342.8809
$ws.Range("I141").Value = 3242.7297
$ws.Range("J141").Value = 4084
$ws.Range("K141").Value = 9728.1891
$ws.Range("L141").Value = 12252
$ws.Range("M141").Value = -4548.1891
$ws.Range("N141").Value = -22612

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 26833
$ws.Range("J24").Value = 26833
$ws.Range("L24").Value = 26833
$ws.Range("N24").Value = -27581
$ws.Range("H32").Value = 1844.03
$ws.Range("I32").Value = 1703.3226
$ws.Range("J32").Value = 3713.4285
$ws.Range("K32").Value = 1703.3226
$ws.Range("L32").Value = 3713.4285
$ws.Range("M32").Value = -1416.3226
$ws.Range("N32").Value = -4287.4285
$ws.Range("H100").Value = 26833
$ws.Range("J100").Value = 26833
$ws.Range("L100").Value = 26833
$ws.Range("N100").Value = -28997
$ws.Range("H110").Value = 3403
$ws.Range("I110").Value = 2580.625
$ws.Range("J110").Value = 4499.5
$ws.Range("K110").Value = 2580.625
$ws.Range("L110").Value = 4499.5
$ws.Range("M110").Value = -535.625
$ws.Range("N110").Value = -8589.5
$ws.Range("H132").Value = 5799.531
$ws.Range("I132").Value = 6122.353
$ws.Range("J132").Value = 5067.8
$ws.Range("K132").Value = 18367.059
$ws.Range("L132").Value = 15203.4
$ws.Range("M132").Value = -15837.059
$ws.Range("N132").Value = -20263.4
$ws.Range("H139").Value = 219572.28
$ws.Range("J139").Value = 219572.28
$ws.Range("L139").Value = 219572.28
$ws.Range("N139").Value = -229852.28

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 39638.332
$ws.Range("J92").Value = 39638.332
$ws.Range("L92").Value = 39638.332
$ws.Range("N92").Value = -44630.332
$ws.Range("H134").Value = 15308.4
$ws.Range("I134").Value = 16148.223
$ws.Range("J134").Value = 7750
$ws.Range("K134").Value = 48444.669
$ws.Range("L134").Value = 23250
$ws.Range("M134").Value = -45909.669
$ws.Range("N134").Value = -28320

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H50").Value = 69000
$ws.Range("J50").Value = 69000
$ws.Range("L50").Value = 69000
$ws.Range("N50").Value = -70250
$ws.Range("H51").Value = 15000
$ws.Range("I51").Value = 15000
$ws.Range("K51").Value = 15000
$ws.Range("M51").Value = -14264
$ws.Range("H58").Value = 2423.4
$ws.Range("I58").Value = 1562
$ws.Range("J58").Value = 4433.3335
$ws.Range("K58").Value = 1562
$ws.Range("L58").Value = 4433.3335
$ws.Range("M58").Value = -1359
$ws.Range("N58").Value = -4839.3335
$ws.Range("H59").Value = 98998
$ws.Range("J59").Value = 98998
$ws.Range("L59").Value = 98998
$ws.Range("N59").Value = -101288
$ws.Range("H61").Value = 15000
$ws.Range("I61").Value = 15000
$ws.Range("K61").Value = 15000
$ws.Range("M61").Value = -14652
$ws.Range("H105").Value = 17437.375
$ws.Range("I105").Value = 25100
$ws.Range("J105").Value = 4666.3335
$ws.Range("K105").Value = 25100
$ws.Range("L105").Value = 4666.3335
$ws.Range("M105").Value = -23353
$ws.Range("N105").Value = -8160.3335
$ws.Range("H107").Value = 9244.308000000001
$ws.Range("I107").Value = 10395.044
$ws.Range("J107").Value = 422
$ws.Range("K107").Value = 10395.044
$ws.Range("L107").Value = 422
$ws.Range("M107").Value = -8475.044
$ws.Range("N107").Value = -4262
$ws.Range("H132").Value = 9731.541999999999
$ws.Range("I132").Value = 1269.4054
$ws.Range("J132").Value = 38195.09
$ws.Range("K132").Value = 3808.2162
$ws.Range("L132").Value = 114585.27
$ws.Range("M132").Value = -1278.2162
$ws.Range("N132").Value = -119645.27
$ws.Range("H134").Value = 1498.3658
$ws.Range("I134").Value = 1418.1034
$ws.Range("J134").Value = 1692.3334
$ws.Range("K134").Value = 4254.3102
$ws.Range("L134").Value = 5077.0002
$ws.Range("M134").Value = -1719.3102
$ws.Range("N134").Value = -10147.0002
$ws.Range("H136").Value = 2423.4
$ws.Range("I136").Value = 1562
$ws.Range("J136").Value = 4433.3335
$ws.Range("K136").Value = 4686
$ws.Range("L136").Value = 13300.0005
$ws.Range("M136").Value = -2136
$ws.Range("N136").Value = -18400.0005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1588.7778
$ws.Range("I51").Value = 483.33334
$ws.Range("J51").Value = 2141.5
$ws.Range("K51").Value = 1450.00002
$ws.Range("L51").Value = 6424.5
$ws.Range("M51").Value = -990.0000199999999
$ws.Range("N51").Value = -7344.5
$ws.Range("H113").Value = 7276.65
$ws.Range("I113").Value = 1175
$ws.Range("J113").Value = 7954.6113
$ws.Range("K113").Value = 3525
$ws.Range("L113").Value = 23863.8339
$ws.Range("M113").Value = -1355
$ws.Range("N113").Value = -28203.8339

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8419.612999999999
$ws.Range("I122").Value = 5040.96
$ws.Range("J122").Value = 22497.334
$ws.Range("K122").Value = 15122.88
$ws.Range("L122").Value = 67492.00199999999
$ws.Range("M122").Value = -12672.88
$ws.Range("N122").Value = -72392.00199999999
$ws.Range("H126").Value = 11792.318
$ws.Range("I126").Value = 16199.111
$ws.Range("K126").Value = 48597.333
$ws.Range("M126").Value = -46127.333
$ws.Range("H132").Value = 2737.4211
$ws.Range("I132").Value = 3063.875
$ws.Range("J132").Value = 996.3333
$ws.Range("K132").Value = 9191.625
$ws.Range("L132").Value = 2988.9999
$ws.Range("M132").Value = -6661.625
$ws.Range("N132").Value = -8048.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2004
$ws.Range("I68").Value = 1700.125
$ws.Range("J68").Value = 2490.2
$ws.Range("K68").Value = 1700.125
$ws.Range("L68").Value = 2490.2
$ws.Range("M68").Value = -951.125
$ws.Range("N68").Value = -3988.2
$ws.Range("H71").Value = 2004
$ws.Range("I71").Value = 1700.125
$ws.Range("J71").Value = 2490.2
$ws.Range("K71").Value = 8500.625
$ws.Range("L71").Value = 12451
$ws.Range("M71").Value = -4756.625
$ws.Range("N71").Value = -19939
$ws.Range("H132").Value = 406094.34
$ws.Range("I132").Value = 574984.7
$ws.Range("J132").Value = 6899
$ws.Range("K132").Value = 1724954.1
$ws.Range("L132").Value = 20697
$ws.Range("M132").Value = -1722424.1
$ws.Range("N132").Value = -25757
